$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.0008583669626518464
$ws.Range("C2").Value = 3099.503889238888
$ws.Range("D2").Value = 337.1190423067083
$ws.Range("E2").Value = 9353990175.932438
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 9353993612.556229
